# Append the "2021年" row (row 11) to Sheet1, mirroring the data added in
# the source diff for 企业单位数.xlsx. The sheet is a transposed table:
# column A holds the year label, row 1 holds the category headers, and each
# data row (2..10 originally, 11 after this edit) holds one year's figures
# across columns B:AQ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A11: year label -------------------------------------------------
# Copy A10's formatting (bold/centered/bordered header style, xf index 1)
# onto A11 so the new row matches the look of every other year row, then
# overwrite the copied value with the new year label.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "2021年"

# --- B11:AQ11: numeric figures for 2021 -------------------------------
# Columns E and V are intentionally left untouched: in the source data
# they hold an explicit empty value (no figure reported for that
# industry in that year) rather than a number, same as rows 2/6-10 (E)
# and rows 9-10 (V) above them.
$ws.Range("B11").Value = 2601
$ws.Range("C11").Value = 825
$ws.Range("D11").Value = 285
$ws.Range("F11").Value = 1326
$ws.Range("G11").Value = 2666
$ws.Range("H11").Value = 200
$ws.Range("I11").Value = 796
$ws.Range("J11").Value = 559
$ws.Range("K11").Value = 43455
$ws.Range("L11").Value = 609
$ws.Range("M11").Value = 103
$ws.Range("N11").Value = 9
$ws.Range("O11").Value = 1478
$ws.Range("P11").Value = 472
$ws.Range("Q11").Value = 34
$ws.Range("R11").Value = 222
$ws.Range("S11").Value = 2845
$ws.Range("T11").Value = 297
$ws.Range("U11").Value = 3222
$ws.Range("W11").Value = 33
$ws.Range("X11").Value = 633
$ws.Range("Y11").Value = 914
$ws.Range("Z11").Value = 3224
$ws.Range("AA11").Value = 932
$ws.Range("AB11").Value = 132
$ws.Range("AC11").Value = 16
$ws.Range("AD11").Value = 1589
$ws.Range("AE11").Value = 1814
$ws.Range("AF11").Value = 4566
$ws.Range("AG11").Value = 3496
$ws.Range("AH11").Value = 734
$ws.Range("AI11").Value = 617
$ws.Range("AJ11").Value = 89
$ws.Range("AK11").Value = 2531
$ws.Range("AL11").Value = 543
$ws.Range("AM11").Value = 1570
$ws.Range("AN11").Value = 50
$ws.Range("AO11").Value = 1090
$ws.Range("AP11").Value = 313
$ws.Range("AQ11").Value = 18
